$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = '61.594.86'
$ws.Range("E2").Value = '  -1.81%  '
$ws.Range("D3").Value = '2.972.64'
$ws.Range("E3").Value = '  -0.93%  '
$ws.Range("E4").Value = '  +0.14%  '
Set-TextValue "D5" '540.36'
$ws.Range("E5").Value = '  -1.98%  '
Set-TextValue "D6" '149.66'
$ws.Range("E6").Value = '  -3.00%  '
Set-TextValue "D7" '1.00'
$ws.Range("E7").Value = '  -0.13%  '
Set-TextValue "D8" '0.563'
$ws.Range("E8").Value = '  +2.19%  '
$ws.Range("D9").Value = '2.984.43'
Set-TextValue "D10" '0.113'
$ws.Range("E10").Value = '  +1.27%  '
Set-TextValue "D11" '6.11'
$ws.Range("E11").Value = '  -4.37%  '
Set-TextValue "D12" '0.366'
$ws.Range("E12").Value = '  +1.37%  '
$ws.Range("D13").Value = '3.497.51'
$ws.Range("E13").Value = '  -0.54%  '
$ws.Range("E14").Value = '  +1.57%  '
$ws.Range("D15").Value = '61.666.62'
$ws.Range("E15").Value = '  -1.88%  '
Set-TextValue "D16" '23.77'
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").Value = '2.989.83'
$ws.Range("E17").Value = '  -0.42%  '
Set-TextValue "D18" '0.0000146'
$ws.Range("E18").Value = '  -1.18%  '
Set-TextValue "D19" '5.13'
$ws.Range("E19").Value = '  +1.99%  '
Set-TextValue "D20" '11.95'
$ws.Range("E20").Value = '  +0.59%  '
Set-TextValue "D21" '374.09'
$ws.Range("E21").Value = '  -3.91%  '
Set-TextValue "D22" '6.68'
$ws.Range("E22").Value = '  +1.13%  '
Set-TextValue "D23" '0.999'
$ws.Range("E23").Value = '  -0.17%  '
Set-TextValue "D24" '65.64'
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").Value = '3.112.74'
$ws.Range("E25").Value = '  -2.23%  '
Set-TextValue "D26" '0.467'
$ws.Range("E26").Value = '  +1.53%  '
Set-TextValue "D27" '0.189'
$ws.Range("E27").Value = '  +1.68%  '
Set-TextValue "D28" '0.998'
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").Value = '0.0₃0909'
$ws.Range("E29").Value = '  -5.13%  '
Set-TextValue "D30" '8.15'
$ws.Range("E30").Value = '  -5.19%  '
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("E32").Value = '  -1.19%  '
Set-TextValue "D33" '20.30'
$ws.Range("E33").Value = '  +0.69%  '
Set-TextValue "D34" '159.80'
$ws.Range("E34").Value = '  -0.73%  '
Set-TextValue "D35" '4.53'
$ws.Range("E35").Value = '  -2.81%  '
Set-TextValue "D36" '5.85'
$ws.Range("E36").Value = '  -1.51%  '
Set-TextValue "D37" '1.05'
$ws.Range("E37").Value = '  -3.29%  '
Set-TextValue "D38" '1.26'
$ws.Range("E38").Value = '  -2.52%  '
Set-TextValue "D39" '1.53'
$ws.Range("E39").Value = '  -2.83%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '2.405.78'
$ws.Range("E40").Value = '  -4.25%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D41" '37.22'
$ws.Range("E41").Value = '  -0.74%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D42" '3.86'
$ws.Range("E42").Value = '  +0.14%  '
Set-TextValue "D43" '21.86'
$ws.Range("E43").Value = '  -2.21%  '
Set-TextValue "D44" '0.669'
$ws.Range("E44").Value = '  +1.64%  '
Set-TextValue "D45" '0.0587'
$ws.Range("E45").Value = '  -1.01%  '
Set-TextValue "D46" '5.18'
$ws.Range("E46").Value = '  +2.16%  '
Set-TextValue "D47" '0.997'
$ws.Range("E47").Value = '  -0.10%  '
Set-TextValue "D48" '0.0243'
$ws.Range("E48").Value = '  -1.21%  '
Set-TextValue "D49" '267.54'
$ws.Range("E49").Value = '  -1.24%  '
Set-TextValue "D50" '0.0946'
$ws.Range("E50").Value = '  +0.93%  '
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue "D51" '10.41'
$ws.Range("E51").Value = '  -0.58%  '
